$wb = $excel.ActiveWorkbook
$origActive = $wb.ActiveSheet
$ws = $wb.Worksheets.Item("BVTStL")
$ws.Range("A1").Value = "(Boolean)"
$ws.Activate() | Out-Null
$ws.Range("A2").Select() | Out-Null
$origActive.Activate() | Out-Null
